# Insert a new daily price record for Mango (Vega Modelo de Temuco) above
# the existing row 389, shifting the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("389:389").Insert()

$ws.Range("A389").Value = 10
$ws.Range("B389").Value = "Vega Modelo de Temuco"
$ws.Range("C389").Value = "La Araucanía"
$ws.Range("D389").Value = 44932
$ws.Range("E389").Value = 9
$ws.Range("F389").Value = "Fruta"
$ws.Range("G389").Value = 100108
$ws.Range("H389").Value = "Tropicales y subtropicales"
$ws.Range("I389").Value = 100108002
$ws.Range("J389").Value = "Mango"
$ws.Range("K389").Value = "Sin especificar"
$ws.Range("L389").Value = "Primera"
$ws.Range("M389").Value = 65
$ws.Range("N389").Value = 8000
$ws.Range("O389").Value = 8000
$ws.Range("P389").Value = 8000
$ws.Range("Q389").Value = "`$/bandeja 4 kilos"
$ws.Range("R389").Value = "Perú"
$ws.Range("S389").Value = 2000
$ws.Range("T389").Value = 4
